$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Constituent element of the network"
$ws.Range("C3").Value = "Burials in the same age tend to have similar goods"
$ws.Range("C4").Value = "Burials of the same sex tend to have similar goods"
$ws.Range("C5").Value = "Burials having ritual practice tend to have similar goods"
$ws.Range("C6").Value = "Burials in the same wealth rank tend to have similar goods "
$ws.Range("C7").Value = "Two burials being connected with a third burial"
$ws.Range("C8").Value = "Burials being connected with multiple partners"
$ws.Range("C9").Value = "Burials physically close to each other tend to have similar goods "

# The longer C9 text now wraps onto a third line, so the row grows taller
# (matches the author's Excel auto re-flow after editing the cell text).
$ws.Rows.Item(9).RowHeight = 51

$ws.Range("C10").Select() | Out-Null
